$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '69.936.36'
$ws.Range('E2').Value = '  +1.09%  '

$ws.Range('D3').Value = '3.506.72'
$ws.Range('E3').Value = '  +0.01%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('D4').ClearFormats()
$ws.Range('E4').Value = '  -0.11%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '606.11'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +4.28%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '172.20'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -0.94%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.616'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  -0.90%  '

$ws.Range('D8').Value = '3.501.76'
$ws.Range('E8').Value = '  +0.09%  '

$ws.Range('E9').Value = '  -0.13%  '

$ws.Range('E10').Value = '  +5.94%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '6.72'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  +0.21%  '

$ws.Range('E12').Value = '  -2.64%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '47.20'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  +0.46%  '

$ws.Range('E14').Value = '  +1.08%  '

$ws.Range('D15').Value = '4.070.55'
$ws.Range('E15').Value = '  +0.10%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '622.31'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  -7.73%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '8.39'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  -3.85%  '

$ws.Range('B18').Value = 'WrappedBTC'
$ws.Range('C18').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D18').Value = '69.813.50'
$ws.Range('E18').Value = '  +0.98%  '

$ws.Range('B19').Value = 'WrappedEther'
$ws.Range('C19').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D19').Value = '3.501.11'
$ws.Range('E19').Value = '  -0.17%  '

$ws.Range('E20').Value = '  -2.18%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '17.29'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  -0.94%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.883'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  -2.23%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '9.83'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  -12.12%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '15.77'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  -2.32%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '96.23'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  -1.73%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '3.84'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  -0.63%  '

$ws.Range('E27').Value = '  -0.02%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.59'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  -2.42%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.15'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  -3.05%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '33.14'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  +0.72%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '8.40'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  -3.69%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.06'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  -4.25%  '

$ws.Range('E33').Value = '  -1.81%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '6.97'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  -4.43%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '564.40'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  -5.44%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '10.74'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  -1.43%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '3.52'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  -2.21%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '57.03'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  -0.51%  '

$ws.Range('E39').Value = '  -3.55%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.999'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  -0.01%  '

$ws.Range('E41').Value = '  +3.94%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.0448'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  +2.04%  '

$ws.Range('B43').Value = 'TheGraph'
$ws.Range('C43').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.326'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  -3.04%  '

$ws.Range('B44').Value = 'Maker'
$ws.Range('C44').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D44').Value = '3.326.72'
$ws.Range('E44').Value = '  -2.68%  '

$ws.Range('B45').Value = 'ThetaToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.98'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  +2.61%  '

$ws.Range('B46').Value = 'PEPE'
$ws.Range('C46').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D46').Value = '0.0₃0708'
$ws.Range('E46').Value = '  +0.07%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '32.97'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  -1.26%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.62'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  +0.86%  '

$ws.Range('E49').Value = '  -3.00%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '134.79'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  +2.52%  '

$ws.Range('E51').Value = '  -1.91%  '
